$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns (D,E) for the new quarter + one more, shifting old D:K to F:M
$ws.Columns("D:E").Insert()

# Copy number formatting/styles from the (now-shifted) F:G columns into the new D:E columns
$ws.Range("F7:G102").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Populate new column D/E values (most recent two quarters) per the updated financial data
$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373
$ws.Range("D8").Value2 = 224600
$ws.Range("E8").Value2 = 224100
$ws.Range("D9").Value2 = 164600
$ws.Range("E9").Value2 = 162200
$ws.Range("D10").Value2 = 60000
$ws.Range("E10").Value2 = 61900
$ws.Range("D12").Value2 = "NA"
$ws.Range("E12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("D14").Value2 = 0
$ws.Range("E14").Value2 = 0
$ws.Range("D15").Value2 = 27600
$ws.Range("E15").Value2 = 23100
$ws.Range("D17").Value2 = 201100
$ws.Range("E17").Value2 = 193100
$ws.Range("D18").Value2 = 23500
$ws.Range("E18").Value2 = 31000
$ws.Range("D20").Value2 = 0
$ws.Range("E20").Value2 = 0
$ws.Range("D21").Value2 = 51100
$ws.Range("E21").Value2 = 54100
$ws.Range("D22").Value2 = 14500
$ws.Range("E22").Value2 = 14900
$ws.Range("D23").Value2 = 9000
$ws.Range("E23").Value2 = 16100
$ws.Range("D24").Value2 = -2600
$ws.Range("E24").Value2 = 400
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("D26").Value2 = 11600
$ws.Range("E26").Value2 = 15700
$ws.Range("D27").Value2 = 11000
$ws.Range("E27").Value2 = 14900
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("D29").Value2 = "NA"
$ws.Range("E29").Value2 = "NA"
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("D32").Value2 = 0
$ws.Range("E32").Value2 = 0
$ws.Range("D33").Value2 = 11000
$ws.Range("E33").Value2 = 14900
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("D35").Value2 = 11000
$ws.Range("E35").Value2 = 14900
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373
$ws.Range("D41").Value2 = 12600
$ws.Range("E41").Value2 = 24200
$ws.Range("D42").Value2 = 0
$ws.Range("E42").Value2 = 0
$ws.Range("D43").Value2 = 51900
$ws.Range("E43").Value2 = 53800
$ws.Range("D44").Value2 = 79000
$ws.Range("E44").Value2 = 82100
$ws.Range("D45").Value2 = 1000
$ws.Range("E45").Value2 = 3700
$ws.Range("D46").Value2 = 144500
$ws.Range("E46").Value2 = 163800
$ws.Range("D47").Value2 = 0
$ws.Range("E47").Value2 = 0
$ws.Range("D48").Value2 = 1245100
$ws.Range("E48").Value2 = 1258200
$ws.Range("D49").Value2 = 229300
$ws.Range("E49").Value2 = 231900
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("D52").Value2 = 200
$ws.Range("E52").Value2 = 200
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("D54").Value2 = 1619100
$ws.Range("E54").Value2 = 1654100
$ws.Range("D57").Value2 = 68800
$ws.Range("E57").Value2 = 90500
$ws.Range("D58").Value2 = 2800
$ws.Range("E58").Value2 = 2700
$ws.Range("D59").Value2 = 19700
$ws.Range("E59").Value2 = 31900
$ws.Range("D60").Value2 = 91300
$ws.Range("E60").Value2 = 125100
$ws.Range("D61").Value2 = 793300
$ws.Range("E61").Value2 = 793300
$ws.Range("D62").Value2 = 127800
$ws.Range("E62").Value2 = 131100
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("D66").Value2 = 1024100
$ws.Range("E66").Value2 = 1061400
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("D70").Value2 = 49300
$ws.Range("E70").Value2 = 39500
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("D72").Value2 = 0
$ws.Range("E72").Value2 = 0
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("D76").Value2 = 545700
$ws.Range("E76").Value2 = 553200
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373
$ws.Range("D81").Value2 = 11000
$ws.Range("E81").Value2 = 14900
$ws.Range("D83").Value2 = 27600
$ws.Range("E83").Value2 = 23100
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("D89").Value2 = 14800
$ws.Range("E89").Value2 = 72600
$ws.Range("D91").Value2 = -16400
$ws.Range("E91").Value2 = -15400
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("D94").Value2 = -16200
$ws.Range("E94").Value2 = -15400
$ws.Range("D96").Value2 = -18900
$ws.Range("E96").Value2 = -18900
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("D100").Value2 = -10200
$ws.Range("E100").Value2 = -45000
$ws.Range("D101").Value2 = 0
$ws.Range("E101").Value2 = 0
$ws.Range("D102").Value2 = -11600
$ws.Range("E102").Value2 = 12200
